# Update the cryptocurrency price/volume listing to the latest scraped
# values (GitHub Actions nightly refresh of cryptos.xlsx).
#
# The "Price" column (D) holds numeric-looking text (e.g. "1.00",
# "0.0000249") that must stay TEXT, exactly like the source file (it was
# authored with inline strings, not numbers). Plain `Range.Value = "1.00"`
# lets Excel's COM layer auto-coerce that into the number 1, dropping the
# trailing zeros / exponent formatting. Forcing the cell to Text format
# before the write, then clearing the format again afterwards, keeps the
# literal string while leaving the cell's style untouched (matches the
# original workbook, which has no explicit style on these cells).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Row 23 / 24 swapped in the rankings: Litecoin <-> PEPE ---
Set-TextValue "B23" "PEPE"
Set-TextValue "C23" "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
Set-TextValue "D23" "0.0000146"
Set-TextValue "E23" "  -4.09%  "

Set-TextValue "B24" "Litecoin"
Set-TextValue "C24" "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
Set-TextValue "D24" "83.21"
Set-TextValue "E24" "  -0.24%  "

# --- Row 38 / 39 swapped in the rankings: Filecoin <-> Mantle ---
Set-TextValue "B38" "Mantle"
Set-TextValue "C38" "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-TextValue "D38" "0.997"
Set-TextValue "E38" "  +0.05%  "

Set-TextValue "B39" "Filecoin"
Set-TextValue "C39" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D39" "5.81"
Set-TextValue "E39" "  +0.74%  "

# --- Remaining Price (D) / Volume(1h) (E) refreshes ---
Set-TextValue "D2" "67.787.75"
Set-TextValue "E2" "  +0.14%  "

Set-TextValue "D3" "3.808.19"
Set-TextValue "E3" "  +0.58%  "

Set-TextValue "D4" "1.00"
Set-TextValue "E4" "  +0.07%  "

Set-TextValue "D5" "603.23"
Set-TextValue "E5" "  +1.27%  "

Set-TextValue "D6" "166.17"
Set-TextValue "E6" "  -0.45%  "

Set-TextValue "E7" "  -0.08%  "

Set-TextValue "D8" "0.519"
Set-TextValue "E8" "  -0.08%  "

Set-TextValue "D9" "0.160"
Set-TextValue "E9" "  +0.35%  "

Set-TextValue "E10" "  +0.93%  "

Set-TextValue "D11" "6.38"
Set-TextValue "E11" "  +1.24%  "

Set-TextValue "D12" "0.0000249"
Set-TextValue "E12" "  -1.20%  "

Set-TextValue "D13" "36.00"
Set-TextValue "E13" "  +0.19%  "

Set-TextValue "D14" "4.446.42"
Set-TextValue "E14" "  +0.54%  "

Set-TextValue "D15" "3.812.73"
Set-TextValue "E15" "  -0.23%  "

Set-TextValue "D16" "67.796.52"
Set-TextValue "E16" "  +0.19%  "

Set-TextValue "D17" "18.38"
Set-TextValue "E17" "  -0.55%  "

Set-TextValue "E18" "  +1.88%  "

Set-TextValue "D19" "7.08"
Set-TextValue "E19" "  +0.72%  "

Set-TextValue "D20" "464.45"
Set-TextValue "E20" "  +1.10%  "

Set-TextValue "D21" "9.83"
Set-TextValue "E21" "  -1.89%  "

Set-TextValue "D22" "0.703"
Set-TextValue "E22" "  +1.05%  "

Set-TextValue "D25" "12.12"
Set-TextValue "E25" "  +0.85%  "

Set-TextValue "E26" "  -0.24%  "

Set-TextValue "D27" "10.02"
Set-TextValue "E27" "  +0.09%  "

Set-TextValue "E28" "  -0.09%  "

Set-TextValue "D29" "3.957.46"
Set-TextValue "E29" "  +0.57%  "

Set-TextValue "E30" "  +0.27%  "

Set-TextValue "D31" "7.44"
Set-TextValue "E31" "  +3.35%  "

Set-TextValue "E32" "  -0.02%  "

Set-TextValue "D33" "29.42"
Set-TextValue "E33" "  -0.64%  "

Set-TextValue "E34" "  +0.07%  "

Set-TextValue "E35" "  -0.11%  "

Set-TextValue "D36" "0.0999"
Set-TextValue "E36" "  -0.10%  "

Set-TextValue "E37" "  -0.09%  "

Set-TextValue "E40" "  -3.27%  "

Set-TextValue "E41" "  -0.04%  "

Set-TextValue "D43" "44.72"
Set-TextValue "E43" "  -2.26%  "

Set-TextValue "D44" "47.76"
Set-TextValue "E44" "  -0.72%  "

Set-TextValue "E45" "  +0.11%  "

Set-TextValue "D46" "28.05"
Set-TextValue "E46" "  +4.80%  "

Set-TextValue "E47" "  +1.56%  "

Set-TextValue "D48" "1.38"
Set-TextValue "E48" "  +11.70%  "

Set-TextValue "D49" "8.36"
Set-TextValue "E49" "  +0.52%  "

Set-TextValue "E50" "  +1.52%  "

Set-TextValue "D51" "389.01"
Set-TextValue "E51" "  -1.05%  "
